$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, matching formatting of existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the data for columns I (I0) and J (IF) across rows 2-39
$data = @(@(1,1),@(8,8),@(7,7),@(1,2),@(7,7),@(4,5),@(8,8),@(1,2),@(7,7),@(9,9),@(9,9),@(9,9),@(8,8),@(1,2),@(6,6),@(8,8),@(8,8),@(6,7),@(7,7),@(8,8),@(7,7),@(6,6),@(6,6),@(6,7),@(6,6),@(8,8),@(8,8),@(7,7),@(6,6),@(7,8),@(2,3),@(8,8),@(8,9),@(6,7),@(7,7),@(6,6),@(6,6),@(5,5))
for ($k = 0; $k -lt $data.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}

Write-Output "Applied I0/IF columns"
